$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 4 new rows before row 141, pushing existing rows 141-176 down to 145-180
$ws.Rows("141:144").Insert()

# Row 141: Lapins / Primera
$ws.Cells.Item(141, 1).Value = 8
$ws.Cells.Item(141, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(141, 3).Value = "Coquimbo"
$ws.Cells.Item(141, 4).Value = 44543
$ws.Cells.Item(141, 5).Value = 4
$ws.Cells.Item(141, 6).Value = "Fruta"
$ws.Cells.Item(141, 7).Value = 100103
$ws.Cells.Item(141, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(141, 9).Value = 100103001
$ws.Cells.Item(141, 10).Value = "Cereza"
$ws.Cells.Item(141, 11).Value = "Lapins"
$ws.Cells.Item(141, 12).Value = "Primera"
$ws.Cells.Item(141, 13).Value = 400
$ws.Cells.Item(141, 14).Value = 11000
$ws.Cells.Item(141, 15).Value = 12000
$ws.Cells.Item(141, 16).Value = 11500
$ws.Cells.Item(141, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(141, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(141, 19).Value = 1150
$ws.Cells.Item(141, 20).Value = 10

# Row 142: Lapins / Segunda
$ws.Cells.Item(142, 1).Value = 8
$ws.Cells.Item(142, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(142, 3).Value = "Coquimbo"
$ws.Cells.Item(142, 4).Value = 44543
$ws.Cells.Item(142, 5).Value = 4
$ws.Cells.Item(142, 6).Value = "Fruta"
$ws.Cells.Item(142, 7).Value = 100103
$ws.Cells.Item(142, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(142, 9).Value = 100103001
$ws.Cells.Item(142, 10).Value = "Cereza"
$ws.Cells.Item(142, 11).Value = "Lapins"
$ws.Cells.Item(142, 12).Value = "Segunda"
$ws.Cells.Item(142, 13).Value = 400
$ws.Cells.Item(142, 14).Value = 9000
$ws.Cells.Item(142, 15).Value = 10000
$ws.Cells.Item(142, 16).Value = 9500
$ws.Cells.Item(142, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(142, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(142, 19).Value = 950
$ws.Cells.Item(142, 20).Value = 10

# Row 143: Rainier / Primera
$ws.Cells.Item(143, 1).Value = 8
$ws.Cells.Item(143, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(143, 3).Value = "Coquimbo"
$ws.Cells.Item(143, 4).Value = 44543
$ws.Cells.Item(143, 5).Value = 4
$ws.Cells.Item(143, 6).Value = "Fruta"
$ws.Cells.Item(143, 7).Value = 100103
$ws.Cells.Item(143, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(143, 9).Value = 100103001
$ws.Cells.Item(143, 10).Value = "Cereza"
$ws.Cells.Item(143, 11).Value = "Rainier"
$ws.Cells.Item(143, 12).Value = "Primera"
$ws.Cells.Item(143, 13).Value = 400
$ws.Cells.Item(143, 14).Value = 17000
$ws.Cells.Item(143, 15).Value = 18000
$ws.Cells.Item(143, 16).Value = 17500
$ws.Cells.Item(143, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(143, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(143, 19).Value = 1750
$ws.Cells.Item(143, 20).Value = 10

# Row 144: Santina / Primera
$ws.Cells.Item(144, 1).Value = 8
$ws.Cells.Item(144, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(144, 3).Value = "Coquimbo"
$ws.Cells.Item(144, 4).Value = 44543
$ws.Cells.Item(144, 5).Value = 4
$ws.Cells.Item(144, 6).Value = "Fruta"
$ws.Cells.Item(144, 7).Value = 100103
$ws.Cells.Item(144, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(144, 9).Value = 100103001
$ws.Cells.Item(144, 10).Value = "Cereza"
$ws.Cells.Item(144, 11).Value = "Santina"
$ws.Cells.Item(144, 12).Value = "Primera"
$ws.Cells.Item(144, 13).Value = 500
$ws.Cells.Item(144, 14).Value = 11000
$ws.Cells.Item(144, 15).Value = 12000
$ws.Cells.Item(144, 16).Value = 11500
$ws.Cells.Item(144, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(144, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(144, 19).Value = 1150
$ws.Cells.Item(144, 20).Value = 10
